$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells in F1:H1
$ws.Range("F1").Value = "id"
$ws.Range("G1").Value = "source_file"
$ws.Range("H1").Value = "text"

# Copy style of existing header (A1) onto the new header cells so formatting matches
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in row 2 data
$ws.Range("A2").Value = "Ruilin"

# B2 holds the text "3" (not a numeric 3), so force text formatting before
# entering the value, then clear the formatting back off so no extra
# number-format styling is left behind on the cell.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"
$ws.Range("B2").ClearFormats()

$ws.Range("C2").Value = "无"
$ws.Range("D2").Value = "SUG"
$ws.Range("E2").Value = "MET"
$ws.Range("F2").Value = "1269f1fb-9c21-42a9-ae5e-c80f92622adc"
$ws.Range("G2").Value = "Bk6qQGWRb_annotated.xlsx"
$ws.Range("H2").Value = "Then how bootstrap dqn extend the idea to deep learning, followed by the noisy net, bbq, shallow UBE and LS-DQN."
